$wb = $excel.ActiveWorkbook

# "workflow" sheet (sheet1.xml) edits: add the "random slopes" row and
# update the random-intercepts example value.
$wsWorkflow = $wb.Worksheets.Item("workflow")

$wsWorkflow.Range("A7").Value = "random slopes"
$wsWorkflow.Range("C7").Value = "list_rand_slopes"
$wsWorkflow.Range("B4").Value = "(1|store)"
$wsWorkflow.Range("B7").Value = "(TV1|store), (TV2|store)"

# fft_terms example value changed from 2 to 0
$wsWorkflow.Range("B3").Value = 0

# Make "workflow" the active sheet/tab with B4 selected
$wsWorkflow.Activate()
$wsWorkflow.Range("B4").Select()
